$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-27 21:18:14"
$ws.Range("E3").Value = "2026-02-27 21:18:16"
$ws.Range("H3").Value = "'41%"
$ws.Range("O3").Value = "4.4 °C"
$ws.Range("E4").Value = "2026-02-27 21:18:18"
$ws.Range("H4").Value = "'89%"
$ws.Range("E5").Value = "2026-02-27 21:18:21"
$ws.Range("H5").Value = "'44%"
$ws.Range("N5").Value = "0.7 °C 20:52 TU"
$ws.Range("O5").Value = "4.8 °C"
$ws.Range("E6").Value = "2026-02-27 21:18:23"
$ws.Range("H6").Value = "'89%"
$ws.Range("E7").Value = "2026-02-27 21:18:26"
$ws.Range("J7").Value = "1024.7 hPa"
$ws.Range("E8").Value = "2026-02-27 21:18:28"
$ws.Range("H8").Value = "'66%"
$ws.Range("J8").Value = "1024.1 hPa"
$ws.Range("N8").Value = "8.0 °C 20:47 TU"
$ws.Range("E9").Value = "2026-02-27 21:18:30"
$ws.Range("E10").Value = "2026-02-27 21:18:33"
$ws.Range("N10").Value = "6.7 °C 20:59 TU"
$ws.Range("O10").Value = "11.0 °C"
$ws.Range("E11").Value = "2026-02-27 21:18:35"
$ws.Range("O11").Value = "8.6 °C"
$ws.Range("E12").Value = "2026-02-27 21:18:38"
$ws.Range("E13").Value = "2026-02-27 21:18:40"
$ws.Range("E14").Value = "2026-02-27 21:18:42"
$ws.Range("O14").Value = "10.6 °C"
$ws.Range("E15").Value = "2026-02-27 21:18:45"
$ws.Range("E16").Value = "2026-02-27 21:18:47"
$ws.Range("N16").Value = "0.1 °C 20:47 TU"
$ws.Range("E17").Value = "2026-02-27 21:18:49"
$ws.Range("N17").Value = "4.9 °C 20:35 TU"
$ws.Range("E18").Value = "2026-02-27 21:18:52"
$ws.Range("E19").Value = "2026-02-27 21:18:54"
$ws.Range("E20").Value = "2026-02-27 21:18:56"
$ws.Range("N20").Value = "-1.3 °C 20:31 TU"
$ws.Range("O20").Value = "3.1 °C"
$ws.Range("E21").Value = "2026-02-27 21:18:59"
$ws.Range("E22").Value = "2026-02-27 21:19:01"
$ws.Range("E23").Value = "2026-02-27 21:19:04"
$ws.Range("H23").Value = "'42%"
$ws.Range("N23").Value = "1.4 °C 20:31 TU"
$ws.Range("E24").Value = "2026-02-27 21:19:06"
$ws.Range("O24").Value = "10.3 °C"
$ws.Range("E25").Value = "2026-02-27 21:19:09"
$ws.Range("N25").Value = "3.0 °C 20:59 TU"
$ws.Range("E26").Value = "2026-02-27 21:19:11"
$ws.Range("H26").Value = "'46%"
$ws.Range("J26").Value = "1021.5 hPa"
$ws.Range("O26").Value = "10.2 °C"
$ws.Range("E27").Value = "2026-02-27 21:19:13"
$ws.Range("N27").Value = "2.6 °C 20:51 TU"
$ws.Range("E28").Value = "2026-02-27 21:19:16"
$ws.Range("E29").Value = "2026-02-27 21:19:18"
$ws.Range("E30").Value = "2026-02-27 21:19:21"
$ws.Range("E31").Value = "2026-02-27 21:19:23"
$ws.Range("E32").Value = "2026-02-27 21:19:26"
$ws.Range("E33").Value = "2026-02-27 21:19:28"
$ws.Range("O33").Value = "8.7 °C"
$ws.Range("E34").Value = "2026-02-27 21:19:31"
$ws.Range("H34").Value = "'48%"
$ws.Range("O34").Value = "4.6 °C"
$ws.Range("E35").Value = "2026-02-27 21:19:33"
$ws.Range("K35").Value = "15.3 MJ/m2"
$ws.Range("O35").Value = "12.0 °C"
$ws.Range("E36").Value = "2026-02-27 21:19:35"
$ws.Range("H36").Value = "'91%"
$ws.Range("J36").Value = "1024.7 hPa"
$ws.Range("E37").Value = "2026-02-27 21:19:38"
$ws.Range("E38").Value = "2026-02-27 21:19:40"
$ws.Range("K38").Value = "12.6 MJ/m2"
$ws.Range("E39").Value = "2026-02-27 21:19:42"
$ws.Range("H39").Value = "'31%"
$ws.Range("E40").Value = "2026-02-27 21:19:45"
$ws.Range("O40").Value = "9.1 °C"
$ws.Range("E41").Value = "2026-02-27 21:19:47"
$ws.Range("E42").Value = "2026-02-27 21:19:49"
$ws.Range("E43").Value = "2026-02-27 21:19:52"
$ws.Range("O43").Value = "9.4 °C"
$ws.Range("E44").Value = "2026-02-27 21:19:54"
$ws.Range("H44").Value = "'60%"
$ws.Range("E45").Value = "2026-02-27 21:19:56"
$ws.Range("H45").Value = "'44%"
$ws.Range("J45").Value = "1021.9 hPa"
$ws.Range("O45").Value = "11.8 °C"
$ws.Range("E46").Value = "2026-02-27 21:19:59"
$ws.Range("J46").Value = "1024.0 hPa"

Write-Host "Applied 88 cell updates"
